{"js": "// Replace each old value with its corresponding new value, matching the\n// target revision (date stamp + all multiplication problems in the table).\nconst replacements = [\n  [\"2025-12-03 Wednesday\", \"2025-12-04 Thursday\"],\n  [\"69\u00d790=\", \"57\u00d782=\"],\n  [\"97\u00d744=\", \"45\u00d799=\"],\n  [\"13\u00d740=\", \"92\u00d788=\"],\n  [\"63\u00d725=\", \"92\u00d759=\"],\n  [\"82\u00d719=\", \"42\u00d747=\"],\n  [\"44\u00d718=\", \"77\u00d776=\"],\n  [\"38\u00d755=\", \"29\u00d722=\"],\n  [\"23\u00d729=\", \"63\u00d784=\"],\n  [\"21\u00d728=\", \"17\u00d768=\"],\n  [\"76\u00d719=\", \"99\u00d711=\"],\n  [\"51\u00d749=\", \"45\u00d713=\"],\n  [\"91\u00d726=\", \"22\u00d761=\"],\n  [\"60\u00d734=\", \"51\u00d745=\"],\n  [\"84\u00d728=\", \"89\u00d797=\"],\n  [\"51\u00d750=\", \"37\u00d796=\"],\n  [\"33\u00d796=\", \"33\u00d721=\"],\n  [\"89\u00d783=\", \"91\u00d790=\"],\n  [\"88\u00d735=\", \"32\u00d756=\"],\n  [\"46\u00d759=\", \"52\u00d738=\"],\n  [\"78\u00d784=\", \"62\u00d721=\"],\n  [\"91\u00d784=\", \"43\u00d747=\"],\n  [\"22\u00d752=\", \"77\u00d765=\"],\n  [\"55\u00d762=\", \"57\u00d732=\"],\n  [\"62\u00d784=\", \"19\u00d790=\"],\n  [\"66\u00d739=\", \"18\u00d791=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply text replacements throughout the document to match the target revision.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-12-03 Wednesday\", \"2025-12-04 Thursday\"),\n    @(\"69\u00d790=\", \"57\u00d782=\"),\n    @(\"97\u00d744=\", \"45\u00d799=\"),\n    @(\"13\u00d740=\", \"92\u00d788=\"),\n    @(\"63\u00d725=\", \"92\u00d759=\"),\n    @(\"82\u00d719=\", \"42\u00d747=\"),\n    @(\"44\u00d718=\", \"77\u00d776=\"),\n    @(\"38\u00d755=\", \"29\u00d722=\"),\n    @(\"23\u00d729=\", \"63\u00d784=\"),\n    @(\"21\u00d728=\", \"17\u00d768=\"),\n    @(\"76\u00d719=\", \"99\u00d711=\"),\n    @(\"51\u00d749=\", \"45\u00d713=\"),\n    @(\"91\u00d726=\", \"22\u00d761=\"),\n    @(\"60\u00d734=\", \"51\u00d745=\"),\n    @(\"84\u00d728=\", \"89\u00d797=\"),\n    @(\"51\u00d750=\", \"37\u00d796=\"),\n    @(\"33\u00d796=\", \"33\u00d721=\"),\n    @(\"89\u00d783=\", \"91\u00d790=\"),\n    @(\"88\u00d735=\", \"32\u00d756=\"),\n    @(\"46\u00d759=\", \"52\u00d738=\"),\n    @(\"78\u00d784=\", \"62\u00d721=\"),\n    @(\"91\u00d784=\", \"43\u00d747=\"),\n    @(\"22\u00d752=\", \"77\u00d765=\"),\n    @(\"55\u00d762=\", \"57\u00d732=\"),\n    @(\"62\u00d784=\", \"19\u00d790=\"),\n    @(\"66\u00d739=\", \"18\u00d791=\"),\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    [void]$find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2)\n}\n"}
